# Update the public-exposure-sites table.
#   Row 2 (was "Caulfield / Metro Train - Frankston line / ... / old")
#     -> replaced with the new Brighton Beach exposure site entry.
#   Row 3 (was "Caulfield / Metro Train - Frankston line / ... / new")
#     -> replaced with the new Cheltenham venue exposure site entry.
#   Two brand-new rows (4 and 5) are appended for the Melbourne restaurant
#     entry, first published with an incorrect end time ("old") and then
#     corrected ("new").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> Brighton Beach
$ws.Range("A2").Value = "Brighton"
$ws.Range("B2").Value = "Brighton Beach"
$ws.Range("C2").Value = "26/12/20 12pm - 1pm"
$ws.Range("D2").Value = "Case attended beach"
$ws.Range("E2").Value = "new"

# Row 3 -> Two Bob Snob, Cheltenham
$ws.Range("A3").Value = "Cheltenham"
$ws.Range("B3").Value = "Two Bob Snob, 256 Charman Road"
$ws.Range("C3").Value = "22/12/2020 1pm - 2pm"
$ws.Range("D3").Value = "Case attended Venue"
$ws.Range("E3").Value = "new"

# Row 4 (new) -> Fonda Mexican, Melbourne, originally-published (superseded) time
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "Fonda Mexican Flinders Lane  31 Flinders Lane Melbourne"
$ws.Range("C4").Value = "29/12/20 6:00pm-7:00pm"
$ws.Range("D4").Value = "Case attended restaurant"
$ws.Range("E4").Value = "old"

# Row 5 (new) -> Fonda Mexican, Melbourne, corrected time
$ws.Range("A5").Value = "Melbourne"
$ws.Range("B5").Value = "Fonda Mexican Flinders Lane  31 Flinders Lane Melbourne"
$ws.Range("C5").Value = "29/12/20 6:00pm-7:30pm"
$ws.Range("D5").Value = "Case attended restaurant"
$ws.Range("E5").Value = "new"

# The author left the cursor on B4 after the edit.
$ws.Range("B4").Select()

# Columns A-D were re-sized (new content changed the natural best-fit widths);
# column E is untouched. ColumnWidth is internally snapped to whole pixels by
# the host the same way real Excel snaps it, so these inputs land on the
# closest achievable stored widths to the authored ones
# (10.2, 46.8, 21.6, 20.7 "characters").
$ws.Columns.Item(1).ColumnWidth = 9.3
$ws.Columns.Item(2).ColumnWidth = 46.0
$ws.Columns.Item(3).ColumnWidth = 20.8
$ws.Columns.Item(4).ColumnWidth = 19.8
